$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 17:31"

# Updated country data (re-sorted by total cases + refreshed stats).
# Each entry: Row number in the sheet, Country name (col A),
# and Vals = Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes (cols B-H)
$countryRows = @(
    @{ Row = 4; Name = "Estados Unidos"; Vals = @(2559422, 6466, 1068949, 1362726, 0, 107, 127747) },
    @{ Row = 5; Name = "Brasil"; Vals = @(1280335, 281, 697526, 526688, 0, 12, 56121) },
    @{ Row = 7; Name = "India"; Vals = @(526113, 16667, 308123, 201931, 0, 370, 16059) },
    @{ Row = 8; Name = "Reino Unido"; Vals = @(310250, 890, 0, 0, 0, 100, 43514) },
    @{ Row = 19; Name = "Francia"; Vals = @(162936, 0, 75649, 57509, 0, 0, 29778) },
    @{ Row = 44; Name = "Republica Dominicana"; Vals = @(30619, 855, 16666, 13235, 0, 6, 718) },
    @{ Row = 45; Name = "Afganistan"; Vals = @(30616, 165, 10674, 19239, 0, 20, 703) },
    @{ Row = 46; Name = "Panama"; Vals = @(29905, 0, 15270, 14060, 0, 0, 575) },
    @{ Row = 57; Name = "Moldavia"; Vals = @(16080, 304, 8963, 6596, 0, 6, 521) },
    @{ Row = 58; Name = "Honduras"; Vals = @(15994, 628, 1678, 13845, 0, 45, 471) },
    @{ Row = 59; Name = "Azerbaiyan"; Vals = @(15890, 521, 8719, 6978, 0, 6, 193) },
    @{ Row = 60; Name = "Ghana"; Vals = @(15834, 0, 11755, 3976, 0, 0, 103) },
    @{ Row = 61; Name = "Guatemala"; Vals = @(15828, 209, 3028, 12128, 0, 49, 672) },
    @{ Row = 81; Name = "Tayikistan"; Vals = @(5799, 52, 4391, 1356, 0, 0, 52) },
    @{ Row = 84; Name = "Etiopia"; Vals = @(5570, 145, 2015, 3461, 0, 5, 94) },
    @{ Row = 95; Name = "Grecia"; Vals = @(3366, 23, 1374, 1801, 0, 0, 191) },
    @{ Row = 153; Name = "Reunion"; Vals = @(520, 3, 472, 46, 0, 0, 2) },
    @{ Row = 201; Name = "Santa Lucia"; Vals = @(19, 0, 19, 0, 0, 0, 0) },
    @{ Row = 202; Name = "Laos"; Vals = @(19, 0, 19, 0, 0, 0, 0) },
    @{ Row = 203; Name = "Fiyi"; Vals = @(18, 0, 18, 0, 0, 0, 0) },
    @{ Row = 204; Name = "Dominica"; Vals = @(18, 0, 18, 0, 0, 0, 0) },
    @{ Row = 208; Name = "Groenlandia"; Vals = @(13, 0, 13, 0, 0, 0, 0) },
    @{ Row = 209; Name = "Islas Malvinas"; Vals = @(13, 0, 13, 0, 0, 0, 0) },
    @{ Row = 212; Name = "Seychelles"; Vals = @(11, 0, 11, 0, 0, 0, 0) },
    @{ Row = 213; Name = "Montserrat"; Vals = @(11, 0, 10, 0, 0, 0, 1) }
)

foreach ($entry in $countryRows) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.Name
    for ($c = 0; $c -lt $entry.Vals.Length; $c++) {
        $ws.Cells.Item($entry.Row, $c + 2).Value = $entry.Vals[$c]
    }
}

